$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.34'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.06%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.74%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.212'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.86%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07700'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.84%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.298'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.30%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.645'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.36%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9151'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.32%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1233'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '10.49%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1820'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.02%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09158'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.17%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04251'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.98%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1050'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.16%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001257'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.44%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005855'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.36%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1,904.52%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.342'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.26%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.322'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '10.62%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.30%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.79%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04024'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.25%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004265'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '3.48%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.00%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02503'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '3.57%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05299'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.36%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007846'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.79%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1315'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.94%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006655'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-5.54%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.61%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007990'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.22%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3068'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.32%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006702'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.20%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.01%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2832'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '799.98%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003101'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-26.17%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.01%'
